$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 height
$ws.Rows.Item(9).RowHeight = 45

# Set values first in the order that controls shared-string index assignment:
# C9 -> 23, D9 -> 24, B9 -> 25
$c9 = $ws.Range("C9")
$c9.Value = "http://192.168.100.19/thaimaiapp/api/mother/locationUpdate"

$d9 = $ws.Range("D9")
$d9.Value = "mPicmeId:1000000000001" + [char]10 + "latitude:12" + [char]10 + "longitude:11"

$b9 = $ws.Range("B9")
$b9.Value = "Location Update"

# Now apply styles in the order that controls cellXfs index assignment:
# D9 (wrapText only) -> new xf 8
$d9.WrapText = $true

# C9 (fill + center/center) -> new xf 9
$c9.Interior.Color = 16777215
$c9.HorizontalAlignment = -4108
$c9.VerticalAlignment = -4108

# B9 reuses existing style 4 (center/center)
$b9.HorizontalAlignment = -4108
$b9.VerticalAlignment = -4108

$ws.Range("C11").Select() | Out-Null
